# Weekly update: insert a new price-report row for Cilantro
# (Macroferia Regional de Talca) at the top of the data block,
# pushing the existing rows 33-77 down to 34-78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 33 - shifts rows 33:77 down to 34:78
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's record
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 44868
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 100112040
$ws.Cells.Item(33, 7).Value = "Cilantro"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 200
$ws.Cells.Item(33, 11).Value = 7000
$ws.Cells.Item(33, 12).Value = 7000
$ws.Cells.Item(33, 13).Value = 7000
$ws.Cells.Item(33, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(33, 15).Value = "Región del Maule"
$ws.Cells.Item(33, 16).Value = 194
$ws.Cells.Item(33, 17).Value = 36
$ws.Cells.Item(33, 18).Value = "Hortaliza"
